$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 6032
$ws.Range("I15").Value = 6032
$ws.Range("K15").Value = 18096
$ws.Range("M15").Value = -17927
$ws.Range("H62").Value = 6500
$ws.Range("H65").Value = 6500
$ws.Range("H74").Value = 9777.182000000001
$ws.Range("I74").Value = 7633.1665
$ws.Range("J74").Value = 12350
$ws.Range("K74").Value = 7633.1665
$ws.Range("L74").Value = 12350
$ws.Range("M74").Value = -6697.1665
$ws.Range("N74").Value = -14222
$ws.Range("H76").Value = 4000
$ws.Range("I76").Value = 4000
$ws.Range("K76").Value = 4000
$ws.Range("M76").Value = -3685
$ws.Range("H77").Value = 9777.182000000001
$ws.Range("I77").Value = 7633.1665
$ws.Range("J77").Value = 12350
$ws.Range("K77").Value = 38165.8325
$ws.Range("L77").Value = 61750
$ws.Range("M77").Value = -33485.8325
$ws.Range("N77").Value = -71110
$ws.Range("H79").Value = 4000
$ws.Range("I79").Value = 4000
$ws.Range("K79").Value = 4000
$ws.Range("M79").Value = -2908
$ws.Range("H137").Value = 4062.6
$ws.Range("I137").Value = 4180.6665
$ws.Range("K137").Value = 12541.9995
$ws.Range("M137").Value = -9991.999500000002
$ws.Range("H141").Value = 3734.5
$ws.Range("I141").Value = 1979.3334
$ws.Range("K141").Value = 5938.0002
$ws.Range("M141").Value = -758.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8896.885
$ws.Range("I32").Value = 7056.826
$ws.Range("J32").Value = 23004
$ws.Range("K32").Value = 7056.826
$ws.Range("L32").Value = 23004
$ws.Range("M32").Value = -6769.826
$ws.Range("N32").Value = -23578
$ws.Range("H122").Value = 2989
$ws.Range("I122").Value = 2989
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8967
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6517
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4908.2036
$ws.Range("I134").Value = 4647.596
$ws.Range("K134").Value = 13942.788
$ws.Range("M134").Value = -11407.788

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 736
$ws.Range("I6").Value = 1062.5
$ws.Range("J6").Value = 246.25
$ws.Range("K6").Value = 1062.5
$ws.Range("L6").Value = 246.25
$ws.Range("M6").Value = -949.5
$ws.Range("N6").Value = -472.25
$ws.Range("H69").Value = 21626.166
$ws.Range("J69").Value = 29919.334
$ws.Range("L69").Value = 29919.334
$ws.Range("N69").Value = -31417.334
$ws.Range("H72").Value = 21626.166
$ws.Range("J72").Value = 29919.334
$ws.Range("L72").Value = 89758.00199999999
$ws.Range("N72").Value = -97246.00199999999
$ws.Range("H86").Value = 11704.5
$ws.Range("I86").Value = 12176.4
$ws.Range("K86").Value = 12176.4
$ws.Range("M86").Value = -11053.4
$ws.Range("H89").Value = 11704.5
$ws.Range("I89").Value = 12176.4
$ws.Range("K89").Value = 60882
$ws.Range("M89").Value = -55266
$ws.Range("H100").Value = 65774.5
$ws.Range("I100").Value = 48999
$ws.Range("J100").Value = 71366.336
$ws.Range("K100").Value = 48999
$ws.Range("L100").Value = 71366.336
$ws.Range("M100").Value = -47917
$ws.Range("N100").Value = -73530.336
$ws.Range("H112").Value = 48531.4
$ws.Range("J112").Value = 48531.4
$ws.Range("L112").Value = 48531.4
$ws.Range("N112").Value = -51485.4
$ws.Range("H122").Value = 4585.3687
$ws.Range("J122").Value = 4254
$ws.Range("L122").Value = 12762
$ws.Range("N122").Value = -17662

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1391.7826
$ws.Range("J132").Value = 1399.4
$ws.Range("L132").Value = 12594.6
$ws.Range("N132").Value = -17654.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 5393.3335
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 450
$ws.Range("I122").Value = 450
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1350
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1100
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 7210.7827
$ws.Range("I132").Value = 6620.125
$ws.Range("J132").Value = 8560.857
$ws.Range("K132").Value = 19860.375
$ws.Range("L132").Value = 25682.571
$ws.Range("M132").Value = -17330.375
$ws.Range("N132").Value = -30742.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7263.3335
$ws.Range("I46").Value = 7152.143
$ws.Range("J46").Value = 7485.7144
$ws.Range("K46").Value = 7152.143
$ws.Range("L46").Value = 7485.7144
$ws.Range("M46").Value = -6964.143
$ws.Range("N46").Value = -7861.7144
$ws.Range("H55").Value = 50000144
$ws.Range("I55").Value = 83333490
$ws.Range("J55").Value = 132.5
$ws.Range("K55").Value = 83333490
$ws.Range("L55").Value = 132.5
$ws.Range("M55").Value = -83333317
$ws.Range("N55").Value = -478.5
$ws.Range("H132").Value = 5794.3423
$ws.Range("I132").Value = 5938.7744
$ws.Range("K132").Value = 17816.3232
$ws.Range("M132").Value = -15286.3232

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5972.25
$ws.Range("I62").Value = 4944.5
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 4944.5
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -4320.5
$ws.Range("N62").Value = -8248
$ws.Range("H65").Value = 5972.25
$ws.Range("I65").Value = 4944.5
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 24722.5
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -21602.5
$ws.Range("N65").Value = -41240
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H99").Value = 32488.5
$ws.Range("I99").Value = 29975
$ws.Range("J99").Value = 33326.332
$ws.Range("K99").Value = 29975
$ws.Range("L99").Value = 33326.332
$ws.Range("M99").Value = -26980
$ws.Range("N99").Value = -39316.332
